$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sprint 106: fill in the Actual output / Result for the last test case row (row 20)
$ws.Range("F20").Value = "It get displayed the successful "
$ws.Range("G20").Value = "Pass"

# Match the formatting used by the row above (borders/wrap-text styles)
$ws.Range("F19:G19").Copy()
$ws.Range("F20:G20").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row grows slightly to fit the new text
$ws.Rows.Item(20).RowHeight = 24.75

# Leave the selection/scroll position where the editor left it
$ws.Range("F20:G20").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 4
